$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.432885646820068
$ws.Range("B1").Value = 2.965213060379028
$ws.Range("C1").Value = 4.723549842834473
$ws.Range("D1").Value = 1.942668199539185
$ws.Range("E1").Value = 1.228165626525879
